# Auto-generated edit script: updates cryptocurrency price/volume data
# (and swaps the Maker/RenderToken rows) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "1.00", "13.19", "0.0826") need the Text number format applied first
# so the literal text is preserved exactly as in the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "42.897.89"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.539.13"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "304.98"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").Value = "98.39"
$ws.Range("E6").Value = "  +6.78%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.547"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "36.86"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "0.0826"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "7.60"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").Value = "2.928.95"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "2.531.26"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "15.13"
$ws.Range("E16").Value = "  +6.60%  "
$ws.Range("D17").Value = "0.867"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "42.895.70"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "71.62"
$ws.Range("D23").Value = "253.76"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "2.93"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "27.75"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +8.33%  "
$ws.Range("D29").Value = "10.22"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "38.47"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "157.15"
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").Value = "19.43"
$ws.Range("E33").Value = "  +13.24%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "0.0798"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "24.78"
$ws.Range("E39").Value = "  +6.19%  "
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  +25.07%  "
$ws.Range("D42").Value = "3.43"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "3.86"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.091.27"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "86.82"
$ws.Range("E47").Value = "  +2.89%  "
$ws.Range("D48").Value = "8.95"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "2.785.49"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "73.71"
$ws.Range("E50").Value = "  +7.10%  "
$ws.Range("E51").Value = "  +1.66%  "
